$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.025.82"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.513.28"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "2.958.93"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "58.978.97"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "2.502.80"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.421"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "0.0₃0766"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.802"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "281.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.605"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0219"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "1.754.77"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.982"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
